# Apply the cryptos.xlsx data refresh captured by the commit diff.
#
# The D (Price) and E (Volume 1h) columns - and, for four coin pairs that
# swapped ranking order, the B (Coin) and C (Link) columns too - are
# rewritten to match the new scrape. All of these cells are plain text in
# the workbook (prices such as "66.270.06" or "1.00" are not valid Excel
# numbers, and the volume cells are percentages wrapped in literal spaces,
# e.g. "  +2.79%  "). Assigning such a string straight to .Value would let
# Excel auto-convert number-looking text into a real number (dropping
# trailing zeros / introducing float noise), so each value is written with
# a leading apostrophe to force text, and the cell style is then reset to
# "Normal" so the quote-prefix formatting does not linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.Value = "'66.270.06"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  +2.79%  "
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.Value = "'2.684.82"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  +2.64%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("E4")
$cell.Value = "'  -0.04%  "
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.Value = "'610.03"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  +2.71%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.Value = "'160.11"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  +4.19%  "
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("E7")
$cell.Value = "'  -0.14%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("D8")
$cell.Value = "'0.593"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  +0.65%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("E9")
$cell.Value = "'  +9.86%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.Value = "'6.03"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  +4.51%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("D11")
$cell.Value = "'0.407"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  +2.82%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("E12")
$cell.Value = "'  +1.66%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.Value = "'0.0000213"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  +24.30%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.Value = "'30.51"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  +6.37%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("D15")
$cell.Value = "'3.166.32"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  +2.51%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.Value = "'66.079.88"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  +2.56%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.Value = "'2.718.65"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  +2.03%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.Value = "'12.71"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  +2.34%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.Value = "'4.90"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  +2.27%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("D20")
$cell.Value = "'362.51"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  +3.36%  "
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("D21")
$cell.Value = "'7.47"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  +4.83%  "
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'  +0.03%  "
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("D23")
$cell.Value = "'70.26"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  +4.06%  "
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("E24")
$cell.Value = "'  +4.36%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("D25")
$cell.Value = "'0.0000107"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  +18.30%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.Value = "'1.66"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  -2.27%  "
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.Value = "'0.174"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  +6.03%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("E28")
$cell.Value = "'  +1.45%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("D29")
$cell.Value = "'8.19"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  -0.90%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Range("E30")
$cell.Value = "'  +7.89%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Range("D31")
$cell.Value = "'542.65"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  +3.34%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Range("D32")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  -0.35%  "
$cell.Style = "Normal"

# Row 33
$cell = $ws.Range("D33")
$cell.Value = "'1.80"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  -0.17%  "
$cell.Style = "Normal"

# Row 34
$cell = $ws.Range("E34")
$cell.Value = "'  +6.15%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("D35")
$cell.Value = "'5.49"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  -4.18%  "
$cell.Style = "Normal"

# Row 36
$cell = $ws.Range("D36")
$cell.Value = "'0.437"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  +3.18%  "
$cell.Style = "Normal"

# Row 37
$cell = $ws.Range("D37")
$cell.Value = "'20.83"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  +3.90%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Range("B38")
$cell.Value = "'Monero"
$cell.Style = "Normal"
$cell = $ws.Range("C38")
$cell.Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Value = "'162.85"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  -0.90%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Range("B39")
$cell.Value = "'Stacks"
$cell.Style = "Normal"
$cell = $ws.Range("C39")
$cell.Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'2.02"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'  +0.38%  "
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  +0.05%  "
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("B41")
$cell.Value = "'USDe"
$cell.Style = "Normal"
$cell = $ws.Range("C41")
$cell.Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  -0.05%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("B42")
$cell.Value = "'Aave"
$cell.Style = "Normal"
$cell = $ws.Range("C42")
$cell.Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'170.30"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +2.98%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.Value = "'42.48"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  +2.45%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.Value = "'4.26"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  +4.74%  "
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("B45")
$cell.Value = "'Hedera"
$cell.Style = "Normal"
$cell = $ws.Range("C45")
$cell.Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'0.0620"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  +3.44%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("B46")
$cell.Value = "'dogwifhat"
$cell.Style = "Normal"
$cell = $ws.Range("C46")
$cell.Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'2.33"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  +5.30%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.Value = "'23.24"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  +0.66%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.Value = "'0.663"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  +3.91%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.Value = "'0.0267"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  +6.86%  "
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.Value = "'20.10"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  +4.66%  "
$cell.Style = "Normal"

# Row 51
$cell = $ws.Range("D51")
$cell.Value = "'0.0990"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  +1.05%  "
$cell.Style = "Normal"
